{"js": "// Replace the date and each \"A\u00d7B=\" equation in the document with the\n// updated values from the commit. Every source string is unique within\n// the document, so a simple exact search-and-replace per pair is safe.\nconst replacements = [\n  [\"2024-05-15 Wednesday\", \"2024-05-16 Thursday\"],\n  [\"695\u00d78=\", \"287\u00d74=\"],\n  [\"270\u00d77=\", \"829\u00d75=\"],\n  [\"538\u00d76=\", \"425\u00d74=\"],\n  [\"931\u00d78=\", \"445\u00d78=\"],\n  [\"901\u00d75=\", \"264\u00d76=\"],\n  [\"873\u00d74=\", \"231\u00d77=\"],\n  [\"673\u00d72=\", \"172\u00d79=\"],\n  [\"454\u00d73=\", \"560\u00d76=\"],\n  [\"630\u00d76=\", \"726\u00d73=\"],\n  [\"960\u00d76=\", \"292\u00d73=\"],\n  [\"652\u00d75=\", \"727\u00d75=\"],\n  [\"844\u00d76=\", \"536\u00d75=\"],\n  [\"630\u00d75=\", \"232\u00d78=\"],\n  [\"710\u00d78=\", \"751\u00d74=\"],\n  [\"910\u00d75=\", \"423\u00d74=\"],\n  [\"155\u00d79=\", \"857\u00d73=\"],\n  [\"131\u00d77=\", \"823\u00d72=\"],\n  [\"602\u00d72=\", \"497\u00d74=\"],\n  [\"835\u00d77=\", \"910\u00d79=\"],\n  [\"392\u00d72=\", \"112\u00d73=\"],\n  [\"187\u00d78=\", \"182\u00d73=\"],\n  [\"268\u00d75=\", \"344\u00d72=\"],\n  [\"627\u00d77=\", \"353\u00d75=\"],\n  [\"607\u00d76=\", \"387\u00d73=\"],\n  [\"179\u00d72=\", \"410\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('2024-05-15 Wednesday', '2024-05-16 Thursday'),\n    @('695\u00d78=', '287\u00d74='),\n    @('270\u00d77=', '829\u00d75='),\n    @('538\u00d76=', '425\u00d74='),\n    @('931\u00d78=', '445\u00d78='),\n    @('901\u00d75=', '264\u00d76='),\n    @('873\u00d74=', '231\u00d77='),\n    @('673\u00d72=', '172\u00d79='),\n    @('454\u00d73=', '560\u00d76='),\n    @('630\u00d76=', '726\u00d73='),\n    @('960\u00d76=', '292\u00d73='),\n    @('652\u00d75=', '727\u00d75='),\n    @('844\u00d76=', '536\u00d75='),\n    @('630\u00d75=', '232\u00d78='),\n    @('710\u00d78=', '751\u00d74='),\n    @('910\u00d75=', '423\u00d74='),\n    @('155\u00d79=', '857\u00d73='),\n    @('131\u00d77=', '823\u00d72='),\n    @('602\u00d72=', '497\u00d74='),\n    @('835\u00d77=', '910\u00d79='),\n    @('392\u00d72=', '112\u00d73='),\n    @('187\u00d78=', '182\u00d73='),\n    @('268\u00d75=', '344\u00d72='),\n    @('627\u00d77=', '353\u00d75='),\n    @('607\u00d76=', '387\u00d73='),\n    @('179\u00d72=', '410\u00d73='),\n)\n\nforeach ($pair in $pairs) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair[0]\n    $find.Replacement.Text = $pair[1]\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
